$d = $word.ActiveDocument

# 1) "7.1 N-Katmanlı Mimari" -> "6.1 Mimari"
[void]$d.Content.Find.Execute("7.1 N-Katmanlı Mimari", $false, $false, $false, $false, $false, `
    $true, 1, $false, "6.1 Mimari", 2)

# 2) "7.2 Tasarım Desenleri" -> "6.2 Tasarım Desenleri"
[void]$d.Content.Find.Execute("7.2 Tasarım Desenleri", $false, $false, $false, $false, $false, `
    $true, 1, $false, "6.2 Tasarım Desenleri", 2)

# 3) "7.3 Konfigürasyon & Logging" -> "6.3 Konfigürasyon & Logging"
[void]$d.Content.Find.Execute("7.3 Konfigürasyon & Logging", $false, $false, $false, $false, $false, `
    $true, 1, $false, "6.3 Konfigürasyon & Logging", 2)

# 4) "8.1 Varlıklar (Özet)" -> "7.1 Varlıklar (Özet)"
[void]$d.Content.Find.Execute("8.1 Varlıklar (Özet)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "7.1 Varlıklar (Özet)", 2)

# 5) First "Dashboard" item (under the "Talep Sahibi" role, right after "Login")
#    gains a leading "s" -> "sDashboard". There are several "Dashboard" paragraphs
#    in the document (one per role), so walk the paragraphs and only touch the
#    first exact match.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Dashboard`r") {
        $para.Range.InsertBefore("s")
        break
    }
}

Write-Output "Done."
